# "fully functional entry deleter": remove the 4th data entry (row 4) from
# the "Lapa1" sheet. Excel shifts the rows below it up, so the rows that
# used to be 5/6/7 (A=5/6/7) become rows 4/5/6, and the former row 7 is gone.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lapa1")
$ws.Rows.Item(4).Delete()
